$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add/update customer measurements (age, height, weight) for the male row (row 2);
# the BMI / body-density / body-fat formulas already in B2:J2 recalculate automatically.
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 1.59
$ws.Range("D2").Value = 55

# Leave the cursor where the user ended up after entering the data.
$ws.Range("H14").Select()
